$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells, copying the existing header's format
# (bold, centered, bordered) so the new columns look like the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row.
for ($row = 2; $row -le 54; $row++) {
    $ws.Cells.Item($row, 30).Value = 84
    $ws.Cells.Item($row, 31).Value = 78
    $ws.Cells.Item($row, 32).Value = 0
}

Write-Host "done"
